$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 1.77
$ws.Range("AA5").Value = 9.25
$ws.Range("AB5").Value = 19
$ws.Range("AC5").Value = 80
$ws.Range("AD5").Value = 7.7
$ws.Range("AE5").Value = 6.9
$ws.Range("AF5").Value = 8.25
$ws.Range("AH5").Value = 10.75
$ws.Range("AI5").Value = 25
$ws.Range("AJ5").Value = 600
$ws.Range("H5").Value = 4.6
$ws.Range("I5").Value = 1.37
$ws.Range("N5").Value = 1.57
$ws.Range("O5").Value = 2.12
$ws.Range("R5").Value = 1.82
$ws.Range("S5").Value = 1.8
$ws.Range("T5").Value = 20
$ws.Range("V5").Value = 23
$ws.Range("X5").Value = 80
$ws.Range("Y5").Value = 70
$ws.Range("Z5").Value = 13.5
$ws.Range("AA7").Value = 14
$ws.Range("AC7").Value = 400
$ws.Range("AF7").Value = 65
$ws.Range("AI7").Value = 350
$ws.Range("H7").Value = 6.1
$ws.Range("I7").Value = 17
$ws.Range("M7").Value = 3.7
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 1.34
$ws.Range("U7").Value = 4.75
$ws.Range("V7").Value = 11.5
$ws.Range("W7").Value = 5.6
$ws.Range("AA8").Value = 7.5
$ws.Range("G8").Value = 3.65
$ws.Range("H8").Value = 3.65
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 8.75
$ws.Range("T8").Value = 12
$ws.Range("U8").Value = 23
$ws.Range("Z8").Value = 8.75
$ws.Range("K15").Value = 9
$ws.Range("N15").Value = 2.15
$ws.Range("O15").Value = 1.67
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10
$ws.Range("N17").Value = 2.08
$ws.Range("O17").Value = 1.73
$ws.Range("AJ18").Value = 451
$ws.Range("N18").Value = 2.03
$ws.Range("O18").Value = 1.78
$ws.Range("AI19").Value = 40
$ws.Range("V19").Value = 8.25
$ws.Range("AA20").Value = 7.4
$ws.Range("AD20").Value = 14
$ws.Range("AF20").Value = 14.5
$ws.Range("AG20").Value = 75
$ws.Range("AH20").Value = 40
$ws.Range("H20").Value = 3.75
$ws.Range("I20").Value = 4.55
$ws.Range("T20").Value = 7.3
$ws.Range("U20").Value = 8
$ws.Range("W20").Value = 12.5
$ws.Range("X20").Value = 13
$ws.Range("Z20").Value = 11.5
$ws.Range("AA21").Value = 6.9
$ws.Range("AB21").Value = 14.5
$ws.Range("AH21").Value = 37
$ws.Range("H21").Value = 3.55
$ws.Range("I21").Value = 3.95
$ws.Range("M21").Value = 3.15
$ws.Range("R21").Value = 1.7
$ws.Range("S21").Value = 1.91
$ws.Range("T21").Value = 7.4
$ws.Range("U21").Value = 8.75
$ws.Range("W21").Value = 15
$ws.Range("Y21").Value = 25
$ws.Range("Z21").Value = 10.75
$ws.Range("AA23").Value = 8.25
$ws.Range("AB23").Value = 19.5
$ws.Range("AC23").Value = 100
$ws.Range("AD23").Value = 6.7
$ws.Range("AE23").Value = 6.6
$ws.Range("AF23").Value = 8.25
$ws.Range("AG23").Value = 9.5
$ws.Range("AH23").Value = 11.75
$ws.Range("AI23").Value = 28
$ws.Range("AJ23").Value = 800
$ws.Range("G23").Value = 6.5
$ws.Range("H23").Value = 4.2
$ws.Range("I23").Value = 1.44
$ws.Range("L23").Value = 1.24
$ws.Range("M23").Value = 3.3
$ws.Range("N23").Value = 1.7
$ws.Range("O23").Value = 1.91
$ws.Range("R23").Value = 1.9
$ws.Range("S23").Value = 1.72
$ws.Range("T23").Value = 16.5
$ws.Range("U23").Value = 40
$ws.Range("V23").Value = 21
$ws.Range("W23").Value = 150
$ws.Range("X23").Value = 75
$ws.Range("Y23").Value = 70
$ws.Range("Z23").Value = 11.25
$ws.Range("AB24").Value = 11.25
$ws.Range("AC24").Value = 45
$ws.Range("AJ24").Value = 300
$ws.Range("T24").Value = 7.5
$ws.Range("U24").Value = 11.75
$ws.Range("AD28").Value = 7.9
$ws.Range("AE28").Value = 15.5
$ws.Range("AF28").Value = 11.5
$ws.Range("AI28").Value = 45
$ws.Range("N28").Value = 2.25
$ws.Range("O28").Value = 1.5
$ws.Range("T28").Value = 6.1
$ws.Range("X28").Value = 22
$ws.Range("Y28").Value = 40
$ws.Range("Z28").Value = 7.1
$ws.Range("AD29").Value = 7.7
$ws.Range("AF29").Value = 11.5
$ws.Range("Q29").Value = 2.18
$ws.Range("T29").Value = 6.2
$ws.Range("U29").Value = 10
$ws.Range("K31").Value = 13
$ws.Range("AG33").Value = 26
$ws.Range("AH33").Value = 19
$ws.Range("G33").Value = 2.75
$ws.Range("I33").Value = 2.5
$ws.Range("N33").Value = 1.65
$ws.Range("O33").Value = 2.2
$ws.Range("T33").Value = 12
$ws.Range("AA34").Value = 10
$ws.Range("AE34").Value = 34
$ws.Range("AF34").Value = 17
$ws.Range("AG34").Value = 51
$ws.Range("G34").Value = 1.48
$ws.Range("I34").Value = 5.5
$ws.Range("G35").Value = 1.92
$ws.Range("AA36").Value = 7.3
$ws.Range("AB36").Value = 14.5
$ws.Range("AC36").Value = 60
$ws.Range("AD36").Value = 12.5
$ws.Range("AF36").Value = 12.5
$ws.Range("AH36").Value = 32
$ws.Range("AI36").Value = 37
$ws.Range("AJ36").Value = 400
$ws.Range("G36").Value = 1.82
$ws.Range("H36").Value = 3.7
$ws.Range("I36").Value = 3.7
$ws.Range("L36").Value = 1.23
$ws.Range("M36").Value = 3.35
$ws.Range("N36").Value = 1.7
$ws.Range("O36").Value = 1.93
$ws.Range("R36").Value = 1.65
$ws.Range("S36").Value = 1.98
$ws.Range("T36").Value = 8
$ws.Range("U36").Value = 9
$ws.Range("V36").Value = 8.25
$ws.Range("W36").Value = 15
$ws.Range("X36").Value = 14
$ws.Range("Y36").Value = 24
$ws.Range("Z36").Value = 12
$ws.Range("AB37").Value = 13
$ws.Range("AC37").Value = 60
$ws.Range("AD37").Value = 9.5
$ws.Range("AE37").Value = 17
$ws.Range("AH37").Value = 26
$ws.Range("AI37").Value = 32
$ws.Range("AJ37").Value = 450
$ws.Range("G37").Value = 2.27
$ws.Range("I37").Value = 3.1
$ws.Range("L37").Value = 1.32
$ws.Range("M37").Value = 2.85
$ws.Range("N37").Value = 1.93
$ws.Range("O37").Value = 1.7
$ws.Range("R37").Value = 1.7
$ws.Range("S37").Value = 1.91
$ws.Range("Z37").Value = 8.75
